$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.200698666666668
$ws.Range("H2").Value = 24.602096
$ws.Range("I2").Value = 0.3423472217473603
$ws.Range("J2").Value = 0.3423472217473603
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.30349533333333
$ws.Range("N2").Value = 30.910486
$ws.Range("O2").Value = 0.376060741358942
$ws.Range("P2").Value = 0.376060741358942
$ws.Range("Q2").Value = 84.49586044207291
$ws.Range("R2").Value = 760.4627439786561
$ws.Range("S2").Value = 0.1287433500124864
$ws.Range("T2").Value = 0.1287433500124864

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.200698666666668
$ws.Range("H3").Value = 24.602096
$ws.Range("I3").Value = 0.3423472217473603
$ws.Range("J3").Value = 0.3423472217473603
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.927949666666668
$ws.Range("N3").Value = 23.783849
$ws.Range("O3").Value = 0.2893572067197239
$ws.Range("P3").Value = 0.2893572067197239
$ws.Range("Q3").Value = 65.0147262608338
$ws.Range("R3").Value = 585.1325363475041
$ws.Range("S3").Value = 0.09906063581307407
$ws.Range("T3").Value = 0.09906063581307409

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.200698666666668
$ws.Range("H4").Value = 24.602096
$ws.Range("I4").Value = 0.3423472217473603
$ws.Range("J4").Value = 0.3423472217473603
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.167042
$ws.Range("N4").Value = 27.501126
$ws.Range("O4").Value = 0.3345820519213342
$ws.Range("P4").Value = 0.3345820519213342
$ws.Range("Q4").Value = 75.17614910667734
$ws.Range("R4").Value = 676.5853419600961
$ws.Range("S4").Value = 0.1145432359217998
$ws.Range("T4").Value = 0.1145432359217998

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.62814333333333
$ws.Range("H5").Value = 31.88443
$ws.Range("I5").Value = 0.4436835799477486
$ws.Range("J5").Value = 0.4436835799477487
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 10.30349533333333
$ws.Range("N5").Value = 30.910486
$ws.Range("O5").Value = 0.376060741358942
$ws.Range("P5").Value = 0.376060741358942
$ws.Range("Q5").Value = 109.5070252369978
$ws.Range("R5").Value = 985.5632271329802
$ws.Range("S5").Value = 0.1668519760039397
$ws.Range("T5").Value = 0.1668519760039398

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.62814333333333
$ws.Range("H6").Value = 31.88443
$ws.Range("I6").Value = 0.4436835799477486
$ws.Range("J6").Value = 0.4436835799477487
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.927949666666668
$ws.Range("N6").Value = 23.783849
$ws.Range("O6").Value = 0.2893572067197239
$ws.Range("P6").Value = 0.2893572067197239
$ws.Range("Q6").Value = 84.25938539678558
$ws.Range("R6").Value = 758.3344685710701
$ws.Range("S6").Value = 0.1283830413610878
$ws.Range("T6").Value = 0.1283830413610878

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.62814333333333
$ws.Range("H7").Value = 31.88443
$ws.Range("I7").Value = 0.4436835799477486
$ws.Range("J7").Value = 0.4436835799477487
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.167042
$ws.Range("N7").Value = 27.501126
$ws.Range("O7").Value = 0.3345820519213342
$ws.Range("P7").Value = 0.3345820519213342
$ws.Range("Q7").Value = 97.42863631868667
$ws.Range("R7").Value = 876.8577268681801
$ws.Range("S7").Value = 0.148448562582721
$ws.Range("T7").Value = 0.1484485625827211

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.125489
$ws.Range("H8").Value = 15.376467
$ws.Range("I8").Value = 0.2139691983048911
$ws.Range("J8").Value = 0.2139691983048911
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.30349533333333
$ws.Range("N8").Value = 30.910486
$ws.Range("O8").Value = 0.376060741358942
$ws.Range("P8").Value = 0.376060741358942
$ws.Range("Q8").Value = 52.81045199255134
$ws.Range("R8").Value = 475.294067932962
$ws.Range("S8").Value = 0.08046541534251579
$ws.Range("T8").Value = 0.08046541534251581

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.125489
$ws.Range("H9").Value = 15.376467
$ws.Range("I9").Value = 0.2139691983048911
$ws.Range("J9").Value = 0.2139691983048911
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.927949666666668
$ws.Range("N9").Value = 23.783849
$ws.Range("O9").Value = 0.2893572067197239
$ws.Range("P9").Value = 0.2893572067197239
$ws.Range("Q9").Value = 40.63461880905367
$ws.Range("R9").Value = 365.711569281483
$ws.Range("S9").Value = 0.06191352954556195
$ws.Range("T9").Value = 0.06191352954556196

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.125489
$ws.Range("H10").Value = 15.376467
$ws.Range("I10").Value = 0.2139691983048911
$ws.Range("J10").Value = 0.2139691983048911
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.167042
$ws.Range("N10").Value = 27.501126
$ws.Range("O10").Value = 0.3345820519213342
$ws.Range("P10").Value = 0.3345820519213342
$ws.Range("Q10").Value = 46.985572933538
$ws.Range("R10").Value = 422.870156401842
$ws.Range("S10").Value = 0.07159025341681331
$ws.Range("T10").Value = 0.07159025341681331

